$wb = $excel.ActiveWorkbook

# --- 1. Text change: "Ready for handoff" -> "In Translation" ---
# The status string is shared across the Overview sheet (columns E/F,
# row 2 - one column per locale) and each per-locale sheet (column C,
# row 2, "Status"). Replace it everywhere it occurs.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# --- 2. Column width changes ---
# With the shorter status text in place, the "Status"-bearing columns no
# longer need to be as wide, so they are narrowed:
#   Overview!E:F      (zh-cn / de-de status columns)  ~17.22 -> ~13.41 chars
#   zh-cn!C / de-de!C (Status column)                 ~17.22 -> ~13.41 chars
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns("E:F").ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns("C:C").ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns("C:C").ColumnWidth = 12.5
